# ticket 171: quitado ID_SERVICIO de plantillas y puesto ID_ORDEN_SERVICIO
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$cellA4 = $ws.Cells.Item(4, 1)

# Rewrite the visible text (keep leading BOM char + new placeholder name)
$bom = [char]0xFEFF
$newText = $bom + "`${ID_ORDEN_SERVICIO}"
$cellA4.Value = $newText

# First run: the lone BOM character keeps its own (tiny, Apple system) font
$run1 = $cellA4.Characters(1, 1)
$run1.Font.Name = ".AppleSystemUIFont"
$run1.Font.Size = 3.9
$run1.Font.ColorIndex = -4105

# Second run: the placeholder text itself, back to the normal Arial font
$run2len = $cellA4.Characters().Text().Length - 1
$run2 = $cellA4.Characters(2, $run2len)
$run2.Font.Name = "Arial"
$run2.Font.Size = 10
$run2.Font.ColorIndex = -4105

# The rest of that header row picked up the same tiny Apple font as a
# whole-cell font (e.g. via a format paint from A4) while keeping their
# original borders/number formats untouched.
$restOfRow = $ws.Range("B4:H4")
$restOfRow.Font.Name = ".AppleSystemUIFont"
$restOfRow.Font.Size = 3.9
$restOfRow.Font.ColorIndex = -4105
